$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '29.351.46'
$ws.Range('D2').Style = 'Normal'
$ws.Range('E2').Value = '  +0.19%  '
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '1.879.06'
$ws.Range('D3').Style = 'Normal'
$ws.Range('E3').Value = '  +0.32%  '
$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '1.000'
$ws.Range('D4').Style = 'Normal'
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '0.7137'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  +0.41%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '242.21'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  -0.14%  '
$ws.Range('E7').Value = '  +0.02%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.08079'
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value = '  +4.48%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.3129'
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  +0.65%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '25.24'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  +0.91%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.08354'
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  -1.37%  '
$ws.Range('B12').Value = 'Polkadot'
$ws.Range('C12').Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '5.248'
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = '  +0.92%  '
$ws.Range('B13').Value = 'WrappedEther'
$ws.Range('C13').Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '1.848.00'
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = '  -1.12%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '0.7190'
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  +1.19%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '91.49'
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  +0.28%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '6.270'
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = '  +4.86%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '0.000008394'
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Value = '  +0.69%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '29.346.65'
$ws.Range('D18').Style = 'Normal'
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '240.79'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  -0.79%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '13.24'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  +0.19%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '2.124.23'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  +0.29%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '0.9998'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  +0.03%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '7.795'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  -0.11%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '0.1593'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  -1.18%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '163.11'
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  +0.18%  '
$ws.Range('E27').Value = '  +0.68%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '18.54'
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  +0.24%  '
$ws.Range('E29').Value = '  -0.31%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '4.424'
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  +0.37%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '4.344'
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  +0.40%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '1.205'
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  -5.10%  '
$ws.Range('E33').Value = '  +2.23%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '1.952'
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  +1.93%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '0.7522'
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  +1.16%  '
$ws.Range('E36').Value = '  +0.52%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '2.700'
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  +0.63%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.01881'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  +1.20%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '1.280.78'
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  +9.69%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '2.740'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  +0.87%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '6.582'
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  +3.66%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '110.36'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  +3.25%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '0.8906'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  +0.54%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '73.06'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  +0.34%  '
$ws.Range('E45').Value = '  +8.37%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '0.9999'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  +0.03%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '2.016.17'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  -0.13%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '1.803'
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  -0.23%  '
$ws.Range('E49').Value = '  +0.11%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '9.470'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  +1.08%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '0.4368'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  +1.69%  '
